$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 461.2
$ws.Range("I19").Value = 355.63635
$ws.Range("J19").Value = 751.5
$ws.Range("K19").Value = 355.63635
$ws.Range("L19").Value = 751.5
$ws.Range("M19").Value = -180.63635
$ws.Range("N19").Value = -1101.5
$ws.Range("H33").Value = 3367861.2
$ws.Range("I33").Value = 1034.7727
$ws.Range("J33").Value = 18181898
$ws.Range("K33").Value = 1034.7727
$ws.Range("L33").Value = 18181898
$ws.Range("M33").Value = -805.7727
$ws.Range("H38").Value = 1201.3158
$ws.Range("I38").Value = 72.083336
$ws.Range("J38").Value = 3137.1428
$ws.Range("K38").Value = 216.250008
$ws.Range("L38").Value = 9411.428400000001
$ws.Range("M38").Value = 155.749992
$ws.Range("N38").Value = -10155.4284
$ws.Range("H128").Value = 15219.444
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 15219.444
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 15219.444
$ws.Range("N128").Value = -25179.444
$ws.Range("H137").Value = 1145.3214
$ws.Range("I137").Value = 1076.0714
$ws.Range("J137").Value = 1353.0714
$ws.Range("K137").Value = 3228.2142
$ws.Range("L137").Value = 4059.2142
$ws.Range("M137").Value = -678.2142000000003

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 21740554
$ws.Range("I74").Value = 1530.6666
$ws.Range("J74").Value = 62501224
$ws.Range("K74").Value = 1530.6666
$ws.Range("L74").Value = 62501224
$ws.Range("M74").Value = -656.6666
$ws.Range("H77").Value = 21740554
$ws.Range("I77").Value = 1530.6666
$ws.Range("J77").Value = 62501224
$ws.Range("K77").Value = 7653.333000000001
$ws.Range("L77").Value = 312506120
$ws.Range("M77").Value = -3285.333000000001
$ws.Range("H97").Value = 1266.6666
$ws.Range("I97").Value = 1266.6666
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1266.6666
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -770.6666
$ws.Range("N97").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1454.2727
$ws.Range("I20").Value = 1393
$ws.Range("J20").Value = 1585.5714
$ws.Range("K20").Value = 1393
$ws.Range("L20").Value = 1585.5714
$ws.Range("M20").Value = -1146
$ws.Range("N20").Value = -2079.5714
$ws.Range("H74").Value = 26938.428
$ws.Range("I74").Value = 2599
$ws.Range("J74").Value = 30995
$ws.Range("K74").Value = 2599
$ws.Range("L74").Value = 30995
$ws.Range("M74").Value = -1663
$ws.Range("N74").Value = -32867
$ws.Range("H77").Value = 26938.428
$ws.Range("I77").Value = 2599
$ws.Range("J77").Value = 30995
$ws.Range("K77").Value = 7797
$ws.Range("L77").Value = 92985
$ws.Range("M77").Value = -3117
$ws.Range("N77").Value = -102345
$ws.Range("H81").Value = 35900
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 35900
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 35900
$ws.Range("N81").Value = -38022
$ws.Range("H84").Value = 35900
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 35900
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 107700
$ws.Range("N84").Value = -118308
$ws.Range("H105").Value = 2439.9
$ws.Range("I105").Value = 2053.5293
$ws.Range("J105").Value = 4629.3335
$ws.Range("K105").Value = 2053.5293
$ws.Range("L105").Value = 4629.3335
$ws.Range("M105").Value = -306.5293000000001
$ws.Range("N105").Value = -8123.3335
$ws.Range("H132").Value = 42000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 42000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 42000
$ws.Range("N132").Value = -52120
$ws.Range("H138").Value = 34725
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 34725
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 34725
$ws.Range("N138").Value = -45005
$ws.Range("H139").Value = 33756.332
$ws.Range("I139").Value = 28709
$ws.Range("J139").Value = 36280
$ws.Range("K139").Value = 28709
$ws.Range("L139").Value = 36280
$ws.Range("M139").Value = -23569
$ws.Range("N139").Value = -46560
$ws.Range("H140").Value = 39106.08
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39106.08
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39106.08
$ws.Range("N140").Value = -49466.08

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2959803.8
$ws.Range("I16").Value = 5918096.5
$ws.Range("J16").Value = 1511
$ws.Range("K16").Value = 5918096.5
$ws.Range("L16").Value = 1511
$ws.Range("M16").Value = -5917809.5
$ws.Range("H31").Value = 9806819
$ws.Range("I31").Value = 1282.4286
$ws.Range("J31").Value = 31256430
$ws.Range("K31").Value = 1282.4286
$ws.Range("L31").Value = 31256430
$ws.Range("M31").Value = -987.4286
$ws.Range("N31").Value = -31257020
$ws.Range("H34").Value = 9806819
$ws.Range("I34").Value = 1282.4286
$ws.Range("J34").Value = 31256430
$ws.Range("K34").Value = 1282.4286
$ws.Range("L34").Value = 31256430
$ws.Range("M34").Value = -1080.4286
$ws.Range("N34").Value = -31256834
$ws.Range("H86").Value = 2299.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2299.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2299.5
$ws.Range("N86").Value = -4545.5
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2299.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2299.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 11497.5
$ws.Range("N89").Value = -22729.5
$ws.Range("M89").ClearContents()
$ws.Range("H105").Value = 25644204
$ws.Range("I105").Value = 37040572
$ws.Range("J105").Value = 2375
$ws.Range("K105").Value = 37040572
$ws.Range("L105").Value = 2375
$ws.Range("M105").Value = -37038825
$ws.Range("N105").Value = -5869
$ws.Range("H107").Value = 271.67648
$ws.Range("I107").Value = 283.14285
$ws.Range("J107").Value = 268.7037
$ws.Range("K107").Value = 283.14285
$ws.Range("L107").Value = 268.7037
$ws.Range("M107").Value = 1636.85715
$ws.Range("N107").Value = -4108.7037
$ws.Range("H113").Value = 2959803.8
$ws.Range("I113").Value = 5918096.5
$ws.Range("J113").Value = 1511
$ws.Range("K113").Value = 5918096.5
$ws.Range("L113").Value = 1511
$ws.Range("M113").Value = -5915926.5
$ws.Range("H134").Value = 12823659
$ws.Range("I134").Value = 16670064
$ws.Range("J134").Value = 2307
$ws.Range("K134").Value = 50010192
$ws.Range("L134").Value = 6921
$ws.Range("M134").Value = -50007657
$ws.Range("N134").Value = -11991

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2400.3333
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 2680.4
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 8041.200000000001
$ws.Range("M25").Value = -2831
$ws.Range("N25").Value = -8379.200000000001
$ws.Range("H30").Value = 2400.3333
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 2680.4
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 8041.200000000001
$ws.Range("M30").Value = -2898
$ws.Range("N30").Value = -8245.200000000001
$ws.Range("H41").Value = 1688.8889
$ws.Range("I41").Value = 1300
$ws.Range("J41").Value = 1737.5
$ws.Range("K41").Value = 3900
$ws.Range("L41").Value = 5212.5
$ws.Range("M41").Value = -3562
$ws.Range("N41").Value = -5888.5
$ws.Range("H61").Value = 488
$ws.Range("I61").Value = 60
$ws.Range("J61").Value = 595
$ws.Range("K61").Value = 180
$ws.Range("L61").Value = 1785
$ws.Range("M61").Value = 35
$ws.Range("N61").Value = -2215
$ws.Range("H64").Value = 2971.4
$ws.Range("I64").Value = 1300
$ws.Range("J64").Value = 3389.25
$ws.Range("K64").Value = 3900
$ws.Range("L64").Value = 10167.75
$ws.Range("M64").Value = -3630
$ws.Range("N64").Value = -10707.75
$ws.Range("H67").Value = 2971.4
$ws.Range("I67").Value = 1300
$ws.Range("J67").Value = 3389.25
$ws.Range("K67").Value = 3900
$ws.Range("L67").Value = 10167.75
$ws.Range("M67").Value = -2964
$ws.Range("N67").Value = -12039.75
$ws.Range("H70").Value = 2426.0833
$ws.Range("I70").Value = 935.5
$ws.Range("J70").Value = 3916.6667
$ws.Range("K70").Value = 2806.5
$ws.Range("L70").Value = 11750.0001
$ws.Range("M70").Value = -2491.5
$ws.Range("H73").Value = 2426.0833
$ws.Range("I73").Value = 935.5
$ws.Range("J73").Value = 3916.6667
$ws.Range("K73").Value = 2806.5
$ws.Range("L73").Value = 11750.0001
$ws.Range("M73").Value = -1714.5
$ws.Range("H105").Value = 6861.9443
$ws.Range("I105").Value = 5526
$ws.Range("J105").Value = 6940.5293
$ws.Range("K105").Value = 16578
$ws.Range("L105").Value = 20821.5879
$ws.Range("M105").Value = -13957
$ws.Range("N105").Value = -26063.5879
$ws.Range("H119").Value = 123318.71
$ws.Range("I119").Value = 1377.4
$ws.Range("J119").Value = 428172
$ws.Range("K119").Value = 4132.200000000001
$ws.Range("L119").Value = 1284516
$ws.Range("M119").Value = 705.7999999999993
$ws.Range("H121").Value = 1041.5714
$ws.Range("I121").Value = 587.1429000000001
$ws.Range("J121").Value = 1268.7858
$ws.Range("K121").Value = 1761.4287
$ws.Range("L121").Value = 3806.3574
$ws.Range("M121").Value = -451.4287000000002
$ws.Range("N121").Value = -6426.357400000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5800
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 5800
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 5800
$ws.Range("N33").Value = -6304
$ws.Range("H40").Value = 12138.462
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 12138.462
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 12138.462
$ws.Range("N40").Value = -12440.462
$ws.Range("H80").Value = 2524.4736
$ws.Range("I80").Value = 2498.75
$ws.Range("J80").Value = 2568.5715
$ws.Range("K80").Value = 2498.75
$ws.Range("L80").Value = 2568.5715
$ws.Range("M80").Value = -1500.75
$ws.Range("N80").Value = -4564.5715
$ws.Range("H83").Value = 2524.4736
$ws.Range("I83").Value = 2498.75
$ws.Range("J83").Value = 2568.5715
$ws.Range("K83").Value = 12493.75
$ws.Range("L83").Value = 12842.8575
$ws.Range("M83").Value = -7501.75
$ws.Range("N83").Value = -22826.8575

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 66669788
$ws.Range("I68").Value = 3150.2
$ws.Range("J68").Value = 200003060
$ws.Range("K68").Value = 3150.2
$ws.Range("L68").Value = 200003060
$ws.Range("M68").Value = -2401.2
$ws.Range("N68").Value = -200004558
$ws.Range("H71").Value = 66669788
$ws.Range("I71").Value = 3150.2
$ws.Range("J71").Value = 200003060
$ws.Range("K71").Value = 15751
$ws.Range("L71").Value = 1000015300
$ws.Range("M71").Value = -12007
$ws.Range("N71").Value = -1000022788
$ws.Range("H82").Value = 65583.06
$ws.Range("I82").Value = 1547.9
$ws.Range("J82").Value = 172308.33
$ws.Range("K82").Value = 1547.9
$ws.Range("L82").Value = 172308.33
$ws.Range("M82").Value = -1186.9
$ws.Range("N82").Value = -173030.33
$ws.Range("H85").Value = 65583.06
$ws.Range("I85").Value = 1547.9
$ws.Range("J85").Value = 172308.33
$ws.Range("K85").Value = 1547.9
$ws.Range("L85").Value = 172308.33
$ws.Range("M85").Value = -299.9000000000001
$ws.Range("N85").Value = -174804.33

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18185116
$ws.Range("I81").Value = 2833.6667
$ws.Range("J81").Value = 25003472
$ws.Range("K81").Value = 5667.3334
$ws.Range("L81").Value = 50006944
$ws.Range("M81").Value = -4606.3334
$ws.Range("N81").Value = -50009066
$ws.Range("H84").Value = 18185116
$ws.Range("I84").Value = 2833.6667
$ws.Range("J84").Value = 25003472
$ws.Range("K84").Value = 28336.667
$ws.Range("L84").Value = 250034720
$ws.Range("M84").Value = -23032.667
$ws.Range("N84").Value = -250045328
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 40400
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40400
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40400
$ws.Range("N138").Value = -50680

